# Parameters.xlsx - add AS/ACS pheromone-update comparison table,
# and highlight the "important" columns (Baterina, Khaluf, Koner, Tian)
# in yellow, per the commit message:
#   "Ook AS/ACS toegevoegd. De gele zijn de artikelen die me belangrijk lijken."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")
$ws.Activate()

# --- Highlight the "important" author columns (B, E, F, H) in yellow ---
# Header row first: these cells already carry the bold header style, so
# this produces the bold+yellow style used for row 1.
$ws.Range("B1").Interior.Color = 65535
$ws.Range("E1").Interior.Color = 65535
$ws.Range("F1").Interior.Color = 65535
$ws.Range("H1").Interior.Color = 65535

# Body cells (rows 2-16) for those same columns: plain (non-bold) + yellow.
$ws.Range("B2:B16").Interior.Color = 65535
$ws.Range("E2:E16").Interior.Color = 65535
$ws.Range("F2:F16").Interior.Color = 65535
$ws.Range("H2:H16").Interior.Color = 65535

# --- New "Pheromone update" section (rows 12-16) ---
# Row 12 is a blank spacer row with the bold column-A style.
$ws.Range("A12").Font.Bold = $true

# Row 13: section header.
$ws.Range("A13").Font.Bold = $true
$ws.Range("A13").Value = "Pheromone update:"

# Row 14: Global update yes/no per author.
$ws.Range("A14").Value = "Global update"
$ws.Range("B14").Value = "yes"
$ws.Range("C14").Value = "yes"
$ws.Range("D14").Value = "no"
$ws.Range("E14").Value = "yes"
$ws.Range("F14").Value = "yes"
$ws.Range("G14").Value = "yes"
$ws.Range("H14").Value = "yes"

# Row 15: Local update yes/no per author.
$ws.Range("A15").Value = "Local update"
$ws.Range("B15").Value = "yes"
$ws.Range("C15").Value = "yes"
$ws.Range("D15").Value = "yes"
$ws.Range("E15").Value = "yes"
$ws.Range("F15").Value = "yes"
$ws.Range("G15").Value = "no"
$ws.Range("H15").Value = "yes"

# Row 16: method (ACS/AS/not indicated) per author.
$ws.Range("A16").Value = "method"
$ws.Range("B16").Value = "ACS"
$ws.Range("C16").Value = "ACS"
$ws.Range("D16").Value = "AS"
$ws.Range("E16").Value = "ACS"
$ws.Range("F16").Value = "not indicated"
$ws.Range("G16").Value = "AS"
$ws.Range("H16").Value = "ACS"

# --- Final selection, matching the saved workbook view ---
[void]$ws.Range("N11").Select()
